# Add 2022-Q1 data
# ------------------------------------------------------------------
# 1) Insert a brand-new worksheet named "2022-Q1" right after "2021-Q4"
#    (and therefore right before "总计"), populate it with the
#    per-fund holdings detail, and
# 2) Insert a new row at the top of the "总计" (totals) sheet's data
#    body summarising the 2022-Q1 quarter.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Part 1: new "2022-Q1" worksheet
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Pull over the existing header / index-column formatting (style id 2:
# bold, bordered, centered) from the "2021-Q4" sheet so we don't create
# redundant new style entries.
$q4.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$q4.Range("A2").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holdings detail rows. Codes / amounts are stored as plain text in
# the source data (e.g. fund codes keep their leading zero), so they are
# entered with a leading apostrophe to force text, matching the original
# inlineStr cells rather than letting Excel reinterpret them as numbers.
$rows = @(
    @{A=0; B="090018"; C="大成新锐产业混合";                     D="125.72"; E="89.64"; F="3.73"; G="4.6894"; H=9},
    @{A=1; B="001300"; C="大成睿景灵活配置混合A";                 D="67.47";  E="89.89"; F="3.35"; G="2.2602"; H=10},
    @{A=2; B="001301"; C="大成睿景灵活配置混合C";                 D="27.47";  E="89.89"; F="3.35"; G="0.9202"; H=10},
    @{A=3; B="002258"; C="大成国企改革灵活配置混合";               D="17.17";  E="89.07"; F="4.40"; G="0.7555"; H=7},
    @{A=4; B="010826"; C="大成产业趋势混合A";                     D="8.93";   E="91.32"; F="3.59"; G="0.3206"; H=10},
    @{A=5; B="005589"; C="长信企业精选两年定期开放灵活配置混合";       D="5.84";   E="79.99"; F="2.69"; G="0.1571"; H=10},
    @{A=6; B="010827"; C="大成产业趋势混合C";                     D="2.01";   E="91.32"; F="3.59"; G="0.0722"; H=10},
    @{A=7; B="001291"; C="摩根士丹利华鑫量化多策略股票";             D="2.36";   E="91.74"; F="1.41"; G="0.0333"; H=10}
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row.A
    $newSheet.Cells.Item($r, 2).Value = "'" + $row.B
    $newSheet.Cells.Item($r, 3).Value = $row.C
    $newSheet.Cells.Item($r, 4).Value = "'" + $row.D
    $newSheet.Cells.Item($r, 5).Value = "'" + $row.E
    $newSheet.Cells.Item($r, 6).Value = "'" + $row.F
    $newSheet.Cells.Item($r, 7).Value = "'" + $row.G
    $newSheet.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# ------------------------------------------------------------------
# Part 2: prepend a "2022-Q1" summary row to the "总计" sheet
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Snapshot the 4 existing data rows (rows 2-5) before shifting them down
# one row to make room for the new "2022-Q1" summary row. Rewriting the
# cells (rather than using Rows.Insert, which drags a copy of the header
# row's formatting down with it) keeps the existing "no explicit style"
# look of columns B-D intact.
$oldVals = @()
for ($r = 2; $r -le 5; $r++) {
    $oldVals += ,@($total.Cells.Item($r, 2).Value2, $total.Cells.Item($r, 3).Value2, $total.Cells.Item($r, 4).Value2)
}

# Extend the bold/bordered index-column style (s=2, shared with the rest
# of column A) down onto the new bottom row before writing into it.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $oldVals.Count; $i++) {
    $destRow = 3 + $i
    $vals = $oldVals[$i]
    $total.Cells.Item($destRow, 1).Value = $i + 1
    $total.Cells.Item($destRow, 2).Value = $vals[0]
    $total.Cells.Item($destRow, 3).Value = $vals[1]
    $total.Cells.Item($destRow, 4).Value = $vals[2]
}

# New row 2: the 2022-Q1 summary
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 9.21
